# Appends the 30/12/2025 16:22 scrape batch (LP1912 + LP1912-215 + 6203-6173)
# to the "horarios-141" tracking workbook, mirroring the scraper's own
# update routine: bump the "Última actualización" timestamp and "Total
# filas" counter on each sheet, then append the new rows it found.

$wb = $excel.ActiveWorkbook

$nuevaActualizacion = "Última actualización: 30/12/2025 16:22:28"

# ---------------------------------------------------------------------
# Sheet "LP1912": columns A..G = Hora_Scrap* / Hora_Scrap / Hora_Llegada /
# Línea / Minutos / Parada / Fecha. New rows appended at 363..385.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2, 1).Value = $nuevaActualizacion
$ws1.Cells.Item(3, 1).Value = "Total filas: 384"

$rows1 = New-Object System.Collections.ArrayList
$rows1.Add(@("", "16:22:17", "16:27", "16_SANTA ANA", 5, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "16:29", "10_OLMOS", 7, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "16:35", "23_HERNANDEZ", 13, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "16:37", "11_ETCHEVERRY", 15, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "16:43", "16_P MOR-SANTA ANA", 21, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "16:48", "15_ABASTO", 26, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "16:56", "17_179 Y 38", 34, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "16:57", "10_OLMOS", 35, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:01", "16_SANTA ANA", 39, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:05", "11_ETCHEVERRY", 43, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:05", "215A_EL PATO", 43, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:05", "23_HERNANDEZ", 43, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:10", "10_OLMOS", 48, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:21", "26_HERNANDEZ", 59, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:24", "84_COLONIA URQUIZA-ESC 49", 62, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:29", "14_ABASTO", 67, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:31", "15_ABASTO", 69, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:33", "27_EL RETIRO", 71, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:38", "17_ROMERO", 76, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:41", "16_SANTA ANA", 79, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:45", "15_ABASTO", 83, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:52", "81_EL PELIGRO", 90, "LP1912", "30/12/2025")) | Out-Null
$rows1.Add(@("", "16:22:17", "17:54", "23_HERNANDEZ", 92, "LP1912", "30/12/2025")) | Out-Null

$startRow1 = 363
for ($i = 0; $i -lt $rows1.Count; $i++) {
    $r = $startRow1 + $i
    $vals = $rows1[$i]
    for ($c = 1; $c -le 7; $c++) {
        $ws1.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215": columns A..G = */ Fecha / Hora_Scrap / Hora_Llegada /
# Línea / Minutos / Parada. New row appended at 30.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2, 1).Value = $nuevaActualizacion
$ws2.Cells.Item(3, 1).Value = "Total filas: 29"

$rows2 = New-Object System.Collections.ArrayList
$rows2.Add(@("", "30/12/2025", "16:22:17", "17:05", "215A_EL PATO", 43, "LP1912")) | Out-Null

$startRow2 = 30
for ($i = 0; $i -lt $rows2.Count; $i++) {
    $r = $startRow2 + $i
    $vals = $rows2[$i]
    for ($c = 1; $c -le 7; $c++) {
        $ws2.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

# ---------------------------------------------------------------------
# Sheet "6203-6173": columns A..G = */ Fecha / Hora_Scrap / Hora_Llegada /
# Línea / Minutos / Parada. New rows appended at 52..53.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = $nuevaActualizacion
$ws3.Cells.Item(3, 1).Value = "Total filas: 52"

$rows3 = New-Object System.Collections.ArrayList
$rows3.Add(@("", "30/12/2025", "16:22:28", "16:53", "215B_LP-P MOR-40 Y 115", 31, "L6173")) | Out-Null
$rows3.Add(@("", "30/12/2025", "16:22:28", "17:26", "215A_LA PLATA", 64, "L6173")) | Out-Null

$startRow3 = 52
for ($i = 0; $i -lt $rows3.Count; $i++) {
    $r = $startRow3 + $i
    $vals = $rows3[$i]
    for ($c = 1; $c -le 7; $c++) {
        $ws3.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
